# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 40
    6  = 517
    7  = 46
    8  = 1963
    9  = 60
    11 = 4147
    13 = 270
    17 = 54
    18 = 2857
    19 = 44
    20 = 401
    22 = 13
    25 = 55
    30 = 279
    31 = 1624
    32 = 230
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
